$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value in E9
$ws.Range("E9").Value = 838212

# Add new row 10 data
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "M2_09 Dryad 2020"
$ws.Range("C10").Value = 9678
$ws.Range("D10").Value = 10725
$ws.Range("E10").Value = 855528
$ws.Range("F10").Value = 9946
$ws.Range("G10").Value = 10046
$ws.Range("H10").Value = 10183

# Apply same style as A2:A9 to A10
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
